$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ecm1-Itgb4 ligand-receptor pair stats for all 4x4 sending/target
# cluster combinations (ECs, FAPs, M2, sCs), following Dr Hou's advice.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T numeric statistics.
$data = @(
    @("ECs", "Ecm1", "Itgb4", "ECs", 3, 1, 7.019285333333333, 21.057856, 0.04775214676711106, 0.04775214676711106, 3, 1, 3.201463666666667, 9.604391, 0.2283977768222221, 0.2283977768222221, 22.47198696063289, 202.247882645696, 0.01090648416009663, 0.01090648416009663),
    @("ECs", "Ecm1", "Itgb4", "FAPs", 3, 1, 7.019285333333333, 21.057856, 0.04775214676711106, 0.04775214676711106, 3, 1, 2.758130333333333, 8.274391, 0.1967696347387151, 0.1967696347387151, 19.36010379618844, 174.240934165696, 0.009396172477353961, 0.009396172477353961),
    @("ECs", "Ecm1", "Itgb4", "M2", 3, 1, 7.019285333333333, 21.057856, 0.04775214676711106, 0.04775214676711106, 2, 0.6666666666666666, 0.4432716666666667, 1.329815, 0.0316237426802848, 0.03162374268028481, 3.111450308515556, 28.00305277664, 0.001510101601794314, 0.001510101601794314),
    @("ECs", "Ecm1", "Itgb4", "sCs", 3, 1, 7.019285333333333, 21.057856, 0.04775214676711106, 0.04775214676711106, 3, 1, 7.614187000000001, 22.842561, 0.5432088457587779, 0.543208845758778, 53.44615113435734, 481.0153602092161, 0.02593938852786616, 0.02593938852786616),
    @("FAPs", "Ecm1", "Itgb4", "ECs", 3, 1, 101.506841, 304.520523, 0.6905502966633176, 0.6905502966633176, 3, 1, 3.201463666666667, 9.604391, 0.2283977768222221, 0.2283977768222221, 324.9704633796104, 2924.734170416493, 0.1577201525418277, 0.1577201525418277),
    @("FAPs", "Ecm1", "Itgb4", "FAPs", 3, 1, 101.506841, 304.520523, 0.6905502966633176, 0.6905502966633176, 3, 1, 2.758130333333333, 8.274391, 0.1967696347387151, 0.1967696347387151, 279.9690972029437, 2519.721874826493, 0.1358793296431524, 0.1358793296431524),
    @("FAPs", "Ecm1", "Itgb4", "M2", 3, 1, 101.506841, 304.520523, 0.6905502966633176, 0.6905502966633176, 2, 0.6666666666666666, 0.4432716666666667, 1.329815, 0.0316237426802848, 0.03162374268028481, 44.99510658813834, 404.955959293245, 0.02183778488947509, 0.02183778488947509),
    @("FAPs", "Ecm1", "Itgb4", "sCs", 3, 1, 101.506841, 304.520523, 0.6905502966633176, 0.6905502966633176, 3, 1, 7.614187000000001, 22.842561, 0.5432088457587779, 0.543208845758778, 772.8920691532671, 6956.028622379405, 0.3751130295888624, 0.3751130295888625),
    @("M2", "Ecm1", "Itgb4", "ECs", 3, 1, 28.868883, 86.606649, 0.1963947998341176, 0.1963947998341176, 3, 1, 3.201463666666667, 9.604391, 0.2283977768222221, 0.2283977768222221, 92.42268002175101, 831.804120195759, 0.04485613566155778, 0.04485613566155778),
    @("M2", "Ecm1", "Itgb4", "FAPs", 3, 1, 28.868883, 86.606649, 0.1963947998341176, 0.1963947998341176, 3, 1, 2.758130333333333, 8.274391, 0.1967696347387151, 0.1967696347387151, 79.62414189175099, 716.617277025759, 0.0386445330279424, 0.0386445330279424),
    @("M2", "Ecm1", "Itgb4", "M2", 3, 1, 28.868883, 86.606649, 0.1963947998341176, 0.1963947998341176, 2, 0.6666666666666666, 0.4432716666666667, 1.329815, 0.0316237426802848, 0.03162374268028481, 12.796757882215, 115.170820939935, 0.006210738613700176, 0.006210738613700178),
    @("M2", "Ecm1", "Itgb4", "sCs", 3, 1, 28.868883, 86.606649, 0.1963947998341176, 0.1963947998341176, 3, 1, 7.614187000000001, 22.842561, 0.5432088457587779, 0.543208845758778, 219.813073643121, 1978.317662788089, 0.1066833925309173, 0.1066833925309173),
    @("sCs", "Ecm1", "Itgb4", "ECs", 3, 1, 9.599121999999999, 28.797366, 0.0653027567354537, 0.0653027567354537, 3, 1, 3.201463666666667, 9.604391, 0.2283977768222221, 0.2283977768222221, 30.73124031490067, 276.581162834106, 0.01491500445874002, 0.01491500445874002),
    @("sCs", "Ecm1", "Itgb4", "FAPs", 3, 1, 9.599121999999999, 28.797366, 0.0653027567354537, 0.0653027567354537, 3, 1, 2.758130333333333, 8.274391, 0.1967696347387151, 0.1967696347387151, 26.47562956156733, 238.280666054106, 0.01284959959026639, 0.01284959959026639),
    @("sCs", "Ecm1", "Itgb4", "M2", 3, 1, 9.599121999999999, 28.797366, 0.0653027567354537, 0.0653027567354537, 2, 0.6666666666666666, 0.4432716666666667, 1.329815, 0.0316237426802848, 0.03162374268028481, 4.255018807476667, 38.29516926729, 0.002065117575315223, 0.002065117575315224),
    @("sCs", "Ecm1", "Itgb4", "sCs", 3, 1, 9.599121999999999, 28.797366, 0.0653027567354537, 0.0653027567354537, 3, 1, 7.614187000000001, 22.842561, 0.5432088457587779, 0.543208845758778, 73.08950994381401, 657.8055894943261, 0.03547303511113207, 0.03547303511113208),
)
$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    for ($col = 1; $col -le $rec.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $rec[$col - 1]
    }
}
